# "Generate Report for Handback"
# Updates the timestamp values recorded in the handback status report.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" for the first row.
$overview.Range("G2").Value = "2016-08-23 17:11:23"

# zh-cn sheet: Correspond Handoff Datetime / Correspond Handback DateTime for the first data row.
$zhcn.Range("H2").Value = "2016-08-23 17:11:17"
$zhcn.Range("K2").Value = "2016-08-23 17:11:35"

# de-de sheet: Correspond Handoff Datetime (same value/shared string as Overview!G2)
# and Correspond Handback DateTime for the first data row.
$dede.Range("H2").Value = "2016-08-23 17:11:23"
$dede.Range("K2").Value = "2016-08-23 17:11:42"
